$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 250, pushing existing rows 250-298 down to 251-299
$ws.Rows(250).Insert()

# Populate the newly inserted row 250 with the new observation
$ws.Range("A250").Value = 8
$ws.Range("B250").Value = "Terminal La Palmera de La Serena"
$ws.Range("C250").Value = "Coquimbo"
$ws.Range("D250").Value = 44644
$ws.Range("E250").Value = 4
$ws.Range("F250").Value = "Fruta"
$ws.Range("G250").Value = 100101
$ws.Range("H250").Value = "Berries"
$ws.Range("I250").Value = 100101007
$ws.Range("J250").Value = "Kiwi"
$ws.Range("K250").Value = "Hayward"
$ws.Range("L250").Value = "Especial"
$ws.Range("M250").Value = 10
$ws.Range("N250").Value = 390000
$ws.Range("O250").Value = 400000
$ws.Range("P250").Value = 395000
$ws.Range("Q250").Value = "`$/bins (450 kilos)"
$ws.Range("R250").Value = "Región de O'Higgins"
$ws.Range("S250").Value = 878
$ws.Range("T250").Value = 450
